$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Keycaps
$ws.Range("A6").Value = "Keycaps"
$ws.Range("C6").Value = "https://amzn.eu/d/jj6j4xa"
$ws.Hyperlinks.Add($ws.Range("C6"), "https://amzn.eu/d/jj6j4xa")
$ws.Range("C6").Style = "Hyperlink"

# Row 7: Switches
$ws.Range("A7").Value = "Switches"
$ws.Range("C7").Value = "https://amzn.eu/d/0HXf95k"
$ws.Hyperlinks.Add($ws.Range("C7"), "https://amzn.eu/d/0HXf95k")
$ws.Range("C7").Style = "Hyperlink"

# Move the active selection down to C8 (matches the post-edit selection in the file)
$ws.Range("C8").Select()
